# edit.ps1 - applies the GOA Groundfish Condition 2021 -> (Oct 2020 update) revisions
# Generated from a verified mapping of old/new text fragments taken from the
# document's own OOXML content (see /tmp/work/gen/build_ps1.py).

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $result = $d.Content.Find.Execute(
        $old, $true, $false, $false, $false, $false,
        $true, 1, $false, $new, 2
    )
    if (-not $result) {
        throw "Find/Replace failed for: $($old.Substring(0, [Math]::Min(60, $old.Length)))"
    }
}

# Last updated date
Replace-Text ': October 2021' ': October 2020'

# Description of Indicator body
Replace-Text ': Residual body condition computed from a long-term average of length-weight-based body condition is an indicator of variability in somatic growth (Brodeur et al., 2004) and represents how heavy a fish is per unit body length. Variability in growth can act as a key indicator of population health and can reflect how populations respond to environmental and other factors (Brosset et al., 2017). Positive residual body condition is interpreted to indicate fish in better condition (heavier per unit length) than those with negative residual body condition indicating poorer condition (lighter per unit length). Overall body condition of fishes likely reflects fish growth which can have implications for their subsequent survival (Paul and Paul, 1999; Boldt and Haldorson, 2004).' ': Length-weight residuals represent how heavy a fish is per unit body length and are an indicator of somatic growth variability (Brodeur et al., 2004). Therefore, length-weight residuals represent an integration of prior prey availability and growth conditions. Positive length-weight residuals indicate better condition (i.e., heavier per unit length) and negative residuals indicate poorer condition (i.e., lighter per unit length). Fish condition calculated as length-weight residuals reflects fish growth trajectories which can have implications for biological productivity due to growth, reproduction, and mortality (Paul and Paul, 1999; Boldt and Haldorson, 2004). In addition, variability in growth and consequent body condition can act as a key indicator of population health reflecting how populations respond to environmental and anthropogenic factors (Brosset et al., 2017).'

# Methods label removal + merge
Replace-Text 'Methods: Paired lengths and weights of individual fishes were examined from the Alaska Fisheries Science Center biennial Resource Assessment and Conservation Engineering (AFSC/RACE) - Groundfish Assessment Program’s (GAP) bottom trawl survey of the Gulf of Alaska (GOA). Analyses focused on walleye pollock (' 'Paired lengths and weights of individual fishes were examined from the Alaska Fisheries Science Center biennial Resource Assessment and Conservation Engineering (AFSC/RACE) - Groundfish Assessment Program’s (GAP) bottom trawl survey of the Gulf of Alaska (GOA). Analyses focused on walleye pollock ('

# GitHub availability sentence
Replace-Text 'is fork length (mm) and a bias correction was applied when predicting weights priori to calculating residuals. Stratum mean residuals were weighted in proportion to stratum biomass and stratum-year combinations with samples sizes <10 were eliminated from indicator calculations although they were included when establishing length-weight relationships. A different slope was estimated for each stratum to account for spatial-temporal variation in growth and bottom trawl survey sampling. Length-weight relationships for 100–250 mm fork length (1–2 year old) walleye pollock were established independent of the adult life history stages caught. Bias-corrected weights-at-length (log scale) were estimated from the model and subtracted from observed weights to compute individual residuals per fish. Length-weight residuals were averaged for each stratum and weighted in proportion to INPFC stratum biomass based on stratified area-swept expansion of summer bottom trawl survey catch per unit effort (CPUE). Average length-weight residuals were compared by stratum and year to evaluate spatial variation in fish condition. As in previous years, confidence intervals for the condition indicator reflect uncertainty based on length-weight residuals, but now better reflect sample sizes and stratum biomasses among years. Confidence intervals do not account for uncertainty in stratum biomass estimates. Combinations of stratum and year with <10 samples were used for length-weight relationships but excluded from indicator calculations. Code used to calculate the condition indicator is available on GitHub at (' 'is fork length (mm) and a bias correction was applied when predicting weights priori to calculating residuals. Stratum mean residuals were weighted in proportion to stratum biomass and stratum-year combinations with samples sizes <10 were eliminated from indicator calculations although they were included when establishing length-weight relationships. A different slope was estimated for each stratum to account for spatial-temporal variation in growth and bottom trawl survey sampling. Length-weight relationships for 100–250 mm fork length (1–2 year old) walleye pollock were established independent of the adult life history stages caught. Bias-corrected weights-at-length (log scale) were estimated from the model and subtracted from observed weights to compute individual residuals per fish. Length-weight residuals were averaged for each stratum and weighted in proportion to INPFC stratum biomass based on stratified area-swept expansion of summer bottom trawl survey catch per unit effort (CPUE). Average length-weight residuals were compared by stratum and year to evaluate spatial variation in fish condition. As in previous years, confidence intervals for the condition indicator reflect uncertainty based on length-weight residuals, but now better reflect sample sizes and stratum biomasses among years. Confidence intervals do not account for uncertainty in stratum biomass estimates. Combinations of stratum and year with <10 samples were used for length-weight relationships but excluded from indicator calculations. Code used to calculate the condition indicator is available at ('

# Status and Trends body
Replace-Text ': Residual body condition varied among survey years for all species considered (Figure 2). Fish condition for all seven species were below average in 2021, but with the same condition or reduction in magnitude for most species in the final year relative to 2019. Residual body condition for pollock, Pacific cod, and arrowtooth flounder remained constant relative to 2019. Southern rock sole residual body condition improved over the last four years, but the final two years remained a constant below average condition. Residual body condition for dusky and northern rockfish also improved, but are still below average. Finally, Pacific ocean perch residual body condition is below average and trending downward in the final four years. Prior to 2015, residual body condition indexes of these GOA species vary from survey to survey, cycling between negative and positive residuals with no clear temporal trends. Residual body condition of 100–250 mm walleye pollock in the GOA is strikingly positive during early years in the time series, but has remained mostly neutral or slightly negative since the early 1990s. Overall, GOA fish condition remains below average.' ': Residual body condition varied among survey years for all species considered (Figure 2). Fish condition indicators for all seven species were below average in 2021, but with the same condition or reduction in magnitude for most species in 2021 relative to 2019. Residual body condition for pollock, Pacific cod, and arrowtooth flounder remained constant relative to 2019. Southern rock sole residual body condition improved over the last four years, but the final two years remained a constant below average condition. Residual body condition for dusky and northern rockfish also improved, but are still below average. Finally, Pacific ocean perch residual body condition is below average and trending downward in the final four years. Prior to 2015, residual body condition indexes of these GOA species vary from survey to survey, cycling between negative and positive residuals with no clear temporal trends. Residual body condition of 100–250 mm walleye pollock in the GOA is strikingly positive during early years in the time series, but has remained mostly neutral or slightly negative since the early 1990s. Overall, GOA fish condition remains below average.'

# Figure 2 caption text
Replace-Text 'Figure 2. Biomass-weighted residual body condition index across survey years (1984-2021) for seven Gulf of Alaska groundfish species collected on the National Marine Fisheries Service (NMFS) Alaska Fisheries Science Center (AFSC) Resource Assessment and Conservation Engineering Groundfish Assessment Program (RACE-GAP) standard summer bottom trawl survey. Filled bars denote weighted length-weight residuals, error bars denote two standard errors.' 'Figure 2. Biomass-weighted residual body condition index across survey years (1984-2021) for seven Gulf of Alaska groundfish species collected on the National Marine Fisheries Service (NMFS) Alaska Fisheries Science Center Resource Assessment and Conservation Engineering (AFSC/RACE) Groundfish Assessment Program (GAP) standard summer bottom trawl survey. Filled bars denote weighted length-weight residuals, error bars denote two standard errors.'

# General patterns paragraph
Replace-Text 'The general patterns of above and below average residual body condition index across recent survey years for the Gulf of Alaska as described above were also apparent in the spatial condition indicators across INPFC strata (Figure 3). The relative contribution of stratum-specific residual body condition to the overall trends (indicated by the height of each colored bar segment) does not demonstrate a clear pattern. Although, for many species, the direction of residual body condition (positive or negative) was synchronous among strata within years. For example, residual body condition for small pollock (100 - 250 mm) in Shumagin and Southeast were positive while other locations trended negative. Residual body condition for southern rock sole in Yakutat and Southeast were also positive, while the rest of the regions trended negative. While Pacific cod residuals trended negative again, residual body condition in the Kodiak strata remained positive. All other fish residual body condition was negative across all strata. Patterns of fish distribution are also apparent in the stratum condition indexes. For example, northern rockfish have primarily been collected from the Shumagin and Chirikof strata in recent surveys.' 'The general patterns of above and below average residual body condition index across recent survey years for the GOA as described above were also apparent in the spatial condition indicators across INPFC strata (Figure 3). The relative contribution of stratum-specific residual body condition to the overall trends (indicated by the height of each colored bar segment) does not demonstrate a clear pattern. Although, for many species, the direction of residual body condition (positive or negative) was synchronous among strata within years. For example, residual body condition for small pollock (100–250 mm) in Shumagin and Southeast were positive while other locations trended negative. Residual body condition for southern rock sole in Yakutat and Southeast were also positive, while the rest of the regions trended negative. While Pacific cod residuals trended negative again, residual body condition in the Kodiak strata remained positive. All other fish residual body condition was negative across all strata. Patterns of fish distribution are also apparent in the stratum condition indexes. For example, northern rockfish have primarily been collected from the Shumagin and Chirikof strata in recent surveys.'

# Figure 3 caption text
Replace-Text 'Figure 3. Residual body condition index for seven Gulf of Alaska groundfish species collected on the National Marine Fisheries Service (NMFS) Alaska Fisheries Science Center (AFSC) Resource Assessment and Conservation Engineering Groundfish Assessment Program (RACE-GAP) standard summer bottom trawl survey (1984–2021) grouped by International North Pacific Fisheries Commission (INPFC) statistical sampling strata.' 'Figure 3. Residual body condition index for seven Gulf of Alaska groundfish species collected on the National Marine Fisheries Service (NMFS) Alaska Fisheries Science Center Resource Assessment and Conservation Engineering (AFSC/RACE) Groundfish Assessment Program (GAP) standard summer bottom trawl survey (1984–2021) grouped by International North Pacific Fisheries Commission (INPFC) statistical sampling strata.'

# Factors causing observed trends body
Replace-Text ': Factors that could affect residual fish body condition presented here include temperature, trawl survey timing, stomach fullness, movement in or out of the survey area, or variable somatic growth. Following an unprecedented warming event from 2014 - 2016 (Bond et al., 2015; Stabeno et al., 2019; Barbeaux et al., 2020), there has been a general trend of warming ocean temperatures in the survey area and sea surface temperature anomaly data continue to reflect temperatures above average historical conditions through 2021 (NOAA 2021); these warmer temperatures could be affecting fish growth conditions in this region. Changing ocean conditions along with normal patterns of movement can cause the proportion of the population resident in the sampling area during the annual bottom trawl survey to vary. Recorded changes attributed to the marine heatwave included species abundances, sizes, growth rates, weight/body condition, reproductive success, and species composition (Suryan et al., 2021). Warmer ocean temperatures can lead to lower energy (leaner) prey, increased metabolic needs of younger fish, and therefore slower growth for juveniles, as observed in Pacific cod (Barbeaux et al., 2020). Despite this evidence, it remains difficult to attribute changes in parameters such as somatic growth and fish condition directly to environmental changes (e.g., Brosset et al., 2017). Additionally, spatial and temporal trends in fish growth over the season become confounded with survey progress since the first length-weight data are generally collected in late May and the bottom trawl survey is conducted throughout the summer months moving from west to east. In addition, spatial variability in residual condition may also reflect local environmental features which can influence growth and prey availability in the areas surveyed (e.g., warm core eddies in the central Gulf of Alaska; Atwood et al., 2010). The fish condition computations presented here begin to, but do not wholly, account for spatio-temporal trends in the data contributed by survey sampling logistics nor do they resolve sources of variability in the underlying populations.' ': Factors that could affect residual fish body condition presented here include temperature, trawl survey timing, stomach fullness, movement in or out of the survey area, or variable somatic growth. Following an unprecedented warming event from 2014–2016 (Bond et al., 2015; Stabeno et al., 2019; Barbeaux et al., 2020), there has been a general trend of warming ocean temperatures in the survey area and sea surface temperature anomaly data continue to reflect temperatures above average historical conditions through 2021 (NOAA 2021); these warmer temperatures could be affecting fish growth conditions in this region. Changing ocean conditions along with normal patterns of movement can cause the proportion of the population resident in the sampling area during the annual bottom trawl survey to vary. Recorded changes attributed to the marine heatwave included species abundances, sizes, growth rates, weight/body condition, reproductive success, and species composition (Suryan et al., 2021). Warmer ocean temperatures can lead to lower energy (leaner) prey, increased metabolic needs of younger fish, and therefore slower growth for juveniles, as observed in Pacific cod (Barbeaux et al., 2020). Additionally, spatial and temporal trends in fish growth over the season become confounded with survey progress since the first length-weight data are generally collected in late May and the bottom trawl survey is conducted throughout the summer months moving from west to east. In addition, spatial variability in residual condition may also reflect local environmental features which can influence growth and prey availability in the areas surveyed (e.g., warm core eddies in the central GOA; Atwood et al., 2010). The fish condition computations presented here begin to, but do not wholly, account for spatio-temporal trends in the data contributed by survey sampling logistics nor do they resolve sources of variability in the underlying populations.'

# Implications body
Replace-Text ': Variations in body condition likely have implications for fish survival. In Prince William Sound, the condition of herring prior to the winter may influence their survival (Paul and Paul, 1999). The condition of Gulf of Alaska groundfish may similarly contribute to survival and recruitment. As future years are added to the time series, the relationship between length-weight residuals and subsequent survival will be examined further. It is important to consider that residual body condition for most species in these analyses was computed for all sizes and sexes combined. Requirements for growth and survivorship differ for different fish life stages and some species have sexually dimorphic or regional growth patterns. It may be more informative to examine life-stage (e.g., early juvenile, subadult, and adult phases) and sex-specific body condition in the future.' ': Variations in body condition likely have implications for fish survival. In Prince William Sound, the condition of herring prior to the winter may influence their survival (Paul and Paul, 1999). The condition of GOA groundfish may similarly contribute to survival and recruitment. As future years are added to the time series, the relationship between length-weight residuals and subsequent survival will be examined further. It is important that residual body condition for most species in these analyses was computed for all sizes and sexes combined. Requirements for growth and survivorship differ for different fish life stages and some species have sexually dimorphic or even regional growth patterns. It may be more informative to examine life-stage (e.g., early juvenile, subadult, and adult phases) and sex-specific body condition in the future.'

# Trend toward lowered body condition paragraph
Replace-Text 'The trend toward lowered body condition for many Gulf of Alaska species over the last 3–4 RACE/AFSC GAP bottom trawl surveys is a potential cause for concern. It could indicate poor overwinter survival or may reflect the influence of locally changing environmental conditions depressing fish growth, local production, or survivorship. Indications are that the Warm Blob (Bond et al., 2015; Stabeno et al., 2019) has been followed by subsequent years with elevated water temperatures (e.g., Barbeaux et al., 2020; NOAA, 2021) which may be related to changes in fish condition in the species examined. As we continue to add years of fish condition to the record and expand on our knowledge of the relationships between condition, growth, production, and survival, we hope to gain more insight into the overall health of fish populations in the Gulf of Alaska.' 'The trend toward lowered body condition for many GOA species over the last three to four RACE/AFSC GAP bottom trawl surveys is a potential cause for concern. It could indicate poor overwinter survival or may reflect the influence of locally changing environmental conditions depressing fish growth, local production, or survivorship. Indications are that the Warm Blob (Bond et al., 2015; Stabeno et al., 2019) has been followed by subsequent years with elevated water temperatures (e.g., Barbeaux et al., 2020; NOAA, 2021) which may be related to changes in fish condition in the species examined. As we continue to add years of fish condition to the record and expand on our knowledge of the relationships between condition, growth, production, and survival, we hope to gain more insight into the overall health of fish populations in the GOA.'

# Research priorities body
Replace-Text ': Efforts are underway to redevelop the groundfish condition indicator for next year’s (2022) ESR, using a spatio-temporal model with spatial random effects (VAST). The change is expected to allow more precise biomass expansion, improve estimates of uncertainty, and better account for spatial-temporal trends in fish growth magnified by our sampling design as well as variation in length-weight samples from bottom trawl surveys due to methodological changes (e.g., transition from sex-and-length stratified to random sampling). In 2022, revised indicators using these spatio-temporal analyses will be presented in a retrospective analysis comparing historical and revised condition indicators. Currently, research is being planned across multiple AFSC programs to explore standardization of statistical methods for calculating condition indicators, and to examine relationships among morphometric condition indicators, bioenergetic indicators, and physiological measures of fish condition.' ': Due to programmatic constraints, we did not transition the groundfish condition indicator to use a spatio-temporal model with spatial random effects (VAST) in 2021. Efforts are underway to redevelop the groundfish condition indicator for next year’s (2022) ESR, using a spatio-temporal model with spatial random effects (VAST; Thorson and Barnett, 2017) and this change should provide more precise biomass expansion, improved estimates of uncertainty, and should better account for spatial-temporal variation in length-weight samples from bottom trawl surveys. Revised indicators will be presented alongside a retrospective analysis to compare the current condition indicator to a VAST-based condition indicator. In addition, research is presently being planned to explore standardization of statistical methods for calculating condition indicators, and to examine relationships among morphometric condition indicators, bioenergetic indicators, and physiological measures of fish condition (Hurst et al., 2021). The Condition Congress Steering Committee provided four recommendations for the future of fish condition research at AFSC: intercalibration of existing condition indices, development of projects to link physiological measurements of condition to demographic outcomes, management-directed research, and standardizing formulation and description of metrics (Hurst et al. 2021). Future research priorities should consider this provided guidance.'

# Inline image alt-text (docPr descr) updates for Figures 2 and 3
# Figure 2 InlineShape AlternativeText
$shp2 = $d.InlineShapes.Item(2)
$shp2.AlternativeText = 'Figure 2. Biomass-weighted residual body condition index across survey years (1984-2021) for seven Gulf of Alaska groundfish species collected on the National Marine Fisheries Service (NMFS) Alaska Fisheries Science Center Resource Assessment and Conservation Engineering (AFSC/RACE) Groundfish Assessment Program (GAP) standard summer bottom trawl survey. Filled bars denote weighted length-weight residuals, error bars denote two standard errors.'

# Figure 3 InlineShape AlternativeText
$shp3 = $d.InlineShapes.Item(3)
$shp3.AlternativeText = 'Figure 3. Residual body condition index for seven Gulf of Alaska groundfish species collected on the National Marine Fisheries Service (NMFS) Alaska Fisheries Science Center Resource Assessment and Conservation Engineering (AFSC/RACE) Groundfish Assessment Program (GAP) standard summer bottom trawl survey (1984–2021) grouped by International North Pacific Fisheries Commission (INPFC) statistical sampling strata.'
